# Create Datagrid for forms
# Replace the Fname/Sname/FullName sample data with a generic
# field1/field2/field3 datagrid layout, and populate the third
# column (previously only partially filled) with x/y/z/q values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "field1"
$ws.Range("B1").Value = "field2"

# Fill the new third column's data rows first (x, y, z, q)
$ws.Range("C2").Value = "x"
$ws.Range("C3").Value = "y"
$ws.Range("C4").Value = "z"
$ws.Range("C5").Value = "q"

# Third column header last
$ws.Range("C1").Value = "field3"

# Selection moves to C1
$ws.Range("C1").Select()
